$d = $word.ActiveDocument

# Find the index of the "Documentation: Microsoft Word" bullet paragraph
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Documentation: Microsoft Word*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Documentation: Microsoft Word' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a new (empty) paragraph right after it; this inherits the
# ListParagraph style + numbering (numId 10) of the preceding bullet.
$target.Range.InsertParagraphAfter()

# Structural edits can invalidate cached paragraph/range handles, so
# re-fetch the freshly created paragraph from the Paragraphs collection
# by its position rather than relying on $target.Next.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Version control: Git/Github"
